# word to pdf conversion fixed
# D2 ("RANGE 1" value on row 2) should hold the text "1-2" (same as E2 "RANGE 2"),
# instead of the numeric value 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from E2 (Text number format) onto D2, then set its value as text "1-2".
$ws.Range("D2").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("D2").Value = "1-2"
